$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'288.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.00%"
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = "'1.48%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.927"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.39%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07330"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.68%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.236"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'23.62%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.712"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.35%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'-1.03%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9016"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.18%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.09194"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'19.03%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1693"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.37%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08172"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.55%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03116"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.63%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09939"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.65%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001498"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.65%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005693"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.79%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.534"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.86%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.065"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.72%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3331"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.43%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'0.03%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.160"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.94%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'-11.94%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04531"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.58%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-0.50%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004164"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'3.90%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'3.99%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003394"
$ws.Range("D27").Style = "Normal"
$ws.Range("D39").Value = "'0.01572"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.62%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04450"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.11%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007433"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'2.43%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009541"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-3.76%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1328"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'1.66%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002291"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'13.35%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009089"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-4.44%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006114"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'2.07%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.01%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.424"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'7.85%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-33.27%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.01%"
$ws.Range("E51").Style = "Normal"
